$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H32").Value = 9997.5
$ws_ALC.Range("J32").Value = 9997.5
$ws_ALC.Range("L32").Value = 9997.5
$ws_ALC.Range("N32").Value = -10649.5
$ws_ALC.Range("H33").Value = 179.33333
$ws_ALC.Range("I33").Value = 179.33333
$ws_ALC.Range("K33").Value = 179.33333
$ws_ALC.Range("M33").Value = 49.66667000000001
$ws_ALC.Range("H39").Value = 30.875
$ws_ALC.Range("I39").Value = 30.875
$ws_ALC.Range("J39").Value = 0
$ws_ALC.Range("K39").Value = 92.625
$ws_ALC.Range("L39").Value = 0
$ws_ALC.Range("M39").ClearContents()
$ws_ALC.Range("N39").Value = 203.375
$ws_ALC.Range("H40").Value = 1187.375
$ws_ALC.Range("I40").Value = 1200
$ws_ALC.Range("K40").Value = 1200
$ws_ALC.Range("M40").Value = -1025
$ws_ALC.Range("H51").Value = 11044.6
$ws_ALC.Range("I51").Value = 9635.143
$ws_ALC.Range("K51").Value = 9635.143
$ws_ALC.Range("M51").Value = -9151.143
$ws_ALC.Range("H80").Value = 3233.3333
$ws_ALC.Range("I80").Value = 0
$ws_ALC.Range("J80").Value = 3233.3333
$ws_ALC.Range("K80").Value = 0
$ws_ALC.Range("L80").ClearContents()
$ws_ALC.Range("M80").Value = 9699.999899999999
$ws_ALC.Range("N80").Value = -11695.9999
$ws_ALC.Range("H83").Value = 3233.3333
$ws_ALC.Range("I83").Value = 0
$ws_ALC.Range("J83").Value = 3233.3333
$ws_ALC.Range("K83").Value = 0
$ws_ALC.Range("L83").ClearContents()
$ws_ALC.Range("M83").Value = 29099.9997
$ws_ALC.Range("N83").Value = -39083.9997
$ws_ALC.Range("H98").Value = 1580.4286
$ws_ALC.Range("I98").Value = 975.9
$ws_ALC.Range("K98").Value = 975.9
$ws_ALC.Range("M98").Value = 522.1
$ws_ALC.Range("H112").Value = 3252.111
$ws_ALC.Range("J112").Value = 3252.111
$ws_ALC.Range("L112").Value = 9756.332999999999
$ws_ALC.Range("N112").Value = -11972.333
$ws_ALC.Range("H113").Value = 7911
$ws_ALC.Range("I113").Value = 12251.167
$ws_ALC.Range("J113").Value = 5543.636
$ws_ALC.Range("K113").Value = 12251.167
$ws_ALC.Range("L113").Value = 5543.636
$ws_ALC.Range("M113").Value = -8997.166999999999
$ws_ALC.Range("N113").Value = -12051.636
$ws_ALC.Range("H122").Value = 1580.4286
$ws_ALC.Range("I122").Value = 975.9
$ws_ALC.Range("K122").Value = 2927.7
$ws_ALC.Range("M122").Value = -477.6999999999998
$ws_ALC.Range("H132").Value = 1173.0834
$ws_ALC.Range("I132").Value = 1173.0834
$ws_ALC.Range("K132").Value = 3519.2502
$ws_ALC.Range("M132").Value = -989.2501999999999

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H23").Value = 10000
$ws_ARM.Range("J23").Value = 10000
$ws_ARM.Range("L23").Value = 10000
$ws_ARM.Range("N23").Value = -10518
$ws_ARM.Range("H37").Value = 6285.5713
$ws_ARM.Range("I37").Value = 3499.75
$ws_ARM.Range("K37").Value = 3499.75
$ws_ARM.Range("M37").Value = -3226.75
$ws_ARM.Range("H61").Value = 1499
$ws_ARM.Range("I61").Value = 1499
$ws_ARM.Range("K61").Value = 1499
$ws_ARM.Range("M61").Value = -1287
$ws_ARM.Range("H110").Value = 659.4
$ws_ARM.Range("I110").Value = 599.6667
$ws_ARM.Range("J110").Value = 749
$ws_ARM.Range("K110").Value = 599.6667
$ws_ARM.Range("L110").Value = 749
$ws_ARM.Range("M110").Value = 1445.3333
$ws_ARM.Range("N110").Value = -4839
$ws_ARM.Range("H122").Value = 1765.6666
$ws_ARM.Range("I122").Value = 1765.6666
$ws_ARM.Range("K122").Value = 5296.9998
$ws_ARM.Range("M122").Value = -2846.9998
$ws_ARM.Range("H132").Value = 979.9231
$ws_ARM.Range("I132").Value = 977.0833
$ws_ARM.Range("J132").Value = 1014
$ws_ARM.Range("K132").Value = 2931.2499
$ws_ARM.Range("L132").Value = 3042
$ws_ARM.Range("M132").Value = -401.2498999999998
$ws_ARM.Range("N132").Value = -8102
$ws_ARM.Range("H136").Value = 1499
$ws_ARM.Range("I136").Value = 1499
$ws_ARM.Range("K136").Value = 4497
$ws_ARM.Range("M136").Value = -1947

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H59").Value = 120000
$ws_BSM.Range("J59").Value = 120000
$ws_BSM.Range("L59").Value = 120000
$ws_BSM.Range("N59").Value = -121694
$ws_BSM.Range("H94").Value = 1075.2941
$ws_BSM.Range("I94").Value = 1177.2142
$ws_BSM.Range("K94").Value = 1177.2142
$ws_BSM.Range("M94").Value = -726.2141999999999
$ws_BSM.Range("H107").Value = 877.8
$ws_BSM.Range("I107").Value = 877.8
$ws_BSM.Range("K107").Value = 877.8
$ws_BSM.Range("M107").Value = 1042.2
$ws_BSM.Range("H134").Value = 481.4
$ws_BSM.Range("I134").Value = 481.4
$ws_BSM.Range("K134").Value = 1444.2
$ws_BSM.Range("M134").Value = 1090.8

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H22").Value = 2350.111
$ws_CRP.Range("I22").Value = 650.25
$ws_CRP.Range("K22").Value = 650.25
$ws_CRP.Range("M22").Value = -300.25
$ws_CRP.Range("H36").Value = 1979.3334
$ws_CRP.Range("I36").Value = 942.5
$ws_CRP.Range("J36").Value = 4053
$ws_CRP.Range("K36").Value = 942.5
$ws_CRP.Range("L36").Value = 4053
$ws_CRP.Range("M36").Value = -554.5
$ws_CRP.Range("N36").Value = -4829
$ws_CRP.Range("H40").Value = 1979.3334
$ws_CRP.Range("I40").Value = 942.5
$ws_CRP.Range("J40").Value = 4053
$ws_CRP.Range("K40").Value = 942.5
$ws_CRP.Range("L40").Value = 4053
$ws_CRP.Range("M40").Value = -782.5
$ws_CRP.Range("N40").Value = -4373
$ws_CRP.Range("H50").Value = 22267.75
$ws_CRP.Range("I50").Value = 9083
$ws_CRP.Range("J50").Value = 26662.666
$ws_CRP.Range("K50").Value = 9083
$ws_CRP.Range("L50").Value = 26662.666
$ws_CRP.Range("M50").Value = -8458
$ws_CRP.Range("N50").Value = -27912.666
$ws_CRP.Range("H107").Value = 690.1
$ws_CRP.Range("I107").Value = 724.6667
$ws_CRP.Range("K107").Value = 724.6667
$ws_CRP.Range("M107").Value = 1195.3333
$ws_CRP.Range("H122").Value = 1333.3334
$ws_CRP.Range("I122").Value = 666.6667
$ws_CRP.Range("K122").Value = 2000.0001
$ws_CRP.Range("M122").Value = 449.9999

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H23").Value = 1019.3333
$ws_CUL.Range("I23").Value = 680.1429000000001
$ws_CUL.Range("J23").Value = 1494.2
$ws_CUL.Range("K23").Value = 2040.4287
$ws_CUL.Range("L23").Value = 4482.6
$ws_CUL.Range("M23").Value = -1805.4287
$ws_CUL.Range("N23").Value = -4952.6
$ws_CUL.Range("H97").Value = 595.6667
$ws_CUL.Range("I97").Value = 562.25
$ws_CUL.Range("J97").Value = 662.5
$ws_CUL.Range("K97").Value = 1686.75
$ws_CUL.Range("L97").Value = 1987.5
$ws_CUL.Range("M97").Value = -1190.75
$ws_CUL.Range("N97").Value = -2979.5

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H57").Value = 16666.666
$ws_GSM.Range("H70").Value = 7000
$ws_GSM.Range("J70").Value = 7000
$ws_GSM.Range("L70").Value = 7000
$ws_GSM.Range("N70").Value = -7540
$ws_GSM.Range("H73").Value = 7000
$ws_GSM.Range("J73").Value = 7000
$ws_GSM.Range("L73").Value = 7000
$ws_GSM.Range("N73").Value = -8872
$ws_GSM.Range("H80").Value = 5452.375
$ws_GSM.Range("I80").Value = 2673
$ws_GSM.Range("J80").Value = 7120
$ws_GSM.Range("K80").Value = 2673
$ws_GSM.Range("L80").Value = 7120
$ws_GSM.Range("M80").Value = -1675
$ws_GSM.Range("N80").Value = -9116
$ws_GSM.Range("H83").Value = 5452.375
$ws_GSM.Range("I83").Value = 2673
$ws_GSM.Range("J83").Value = 7120
$ws_GSM.Range("K83").Value = 13365
$ws_GSM.Range("L83").Value = 35600
$ws_GSM.Range("M83").Value = -8373
$ws_GSM.Range("N83").Value = -45584
$ws_GSM.Range("H102").Value = 1168.5555
$ws_GSM.Range("I102").Value = 1190.625
$ws_GSM.Range("K102").Value = 1190.625
$ws_GSM.Range("M102").Value = 431.375
$ws_GSM.Range("H132").Value = 4090.6667
$ws_GSM.Range("I132").Value = 4090.6667
$ws_GSM.Range("K132").Value = 12272.0001
$ws_GSM.Range("M132").Value = -9742.000100000001

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 3919.6924
$ws_LTW.Range("I40").Value = 3913
$ws_LTW.Range("K40").Value = 3913
$ws_LTW.Range("M40").Value = -3777
$ws_LTW.Range("H46").Value = 3693.4
$ws_LTW.Range("I46").Value = 3693.4
$ws_LTW.Range("J46").Value = 0
$ws_LTW.Range("K46").Value = 3693.4
$ws_LTW.Range("L46").Value = 0
$ws_LTW.Range("M46").ClearContents()
$ws_LTW.Range("N46").Value = -3505.4
$ws_LTW.Range("H132").Value = 6829
$ws_LTW.Range("I132").Value = 3100
$ws_LTW.Range("J132").Value = 10558
$ws_LTW.Range("K132").Value = 9300
$ws_LTW.Range("L132").Value = 31674
$ws_LTW.Range("M132").Value = -6770
$ws_LTW.Range("N132").Value = -36734
$ws_LTW.Range("H136").Value = 6250
$ws_LTW.Range("I136").Value = 10000
$ws_LTW.Range("K136").Value = 30000
$ws_LTW.Range("M136").Value = -27450

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H82").Value = 20301
$ws_WVR.Range("J82").Value = 20301
$ws_WVR.Range("L82").Value = 20301
$ws_WVR.Range("N82").Value = -21067
$ws_WVR.Range("H85").Value = 20301
$ws_WVR.Range("J85").Value = 20301
$ws_WVR.Range("L85").Value = 20301
$ws_WVR.Range("N85").Value = -22953
$ws_WVR.Range("H107").Value = 599
$ws_WVR.Range("I107").Value = 599
$ws_WVR.Range("J107").Value = 0
$ws_WVR.Range("K107").Value = 1797
$ws_WVR.Range("L107").ClearContents()
$ws_WVR.Range("N107").Value = 0
$ws_WVR.Range("M107").Value = 123
$ws_WVR.Range("H132").Value = 3428.5715
$ws_WVR.Range("I132").Value = 3166.6667
$ws_WVR.Range("K132").Value = 9500.000100000001
$ws_WVR.Range("M132").Value = -6970.000100000001
